$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (target XML widths: C=54, D=51, G=16, H=33) ---
# Engine's ColumnWidth (char units) = stored xml width - 0.83
$ws.Columns.Item(3).ColumnWidth = 53.17
$ws.Columns.Item(4).ColumnWidth = 50.17
$ws.Columns.Item(7).ColumnWidth = 15.17
$ws.Columns.Item(8).ColumnWidth = 32.17

# --- New fill style (yellow) used to highlight PREMIUM = "Yes" rows ---
# Apply to the E-column cells that become "Yes" further below; this claims
# cellXfs index 3 (fillId=3) before anything else touches the style table.
$ws.Range("E2").Interior.Color = 65535
$ws.Range("E3").Interior.Color = 65535
$ws.Range("E5").Interior.Color = 65535
$ws.Range("E7").Interior.Color = 65535

# --- Keep the OPPORTUNITY ID column stored as text (matches scrape format) ---
# Force text number-format on column A for the data rows, write the values,
# then reset the visible style to Normal so no stray "s" attribute remains.
$ws.Range("A2:A8").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "1331527"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1331527"
$ws.Range("C2").Value = "[CC] Project & Process Management Support"
$ws.Range("D2").Value = "Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "8 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "DHL Group"

# Row 3
$ws.Range("A3").Value = "1331515"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1331515"
$ws.Range("C3").Value = "[SNL] People Analytics"
$ws.Range("D3").Value = "Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany"
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "9 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "DHL Group"

# Row 4
$ws.Range("A4").Value = "1331514"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1331514"
$ws.Range("C4").Value = "Sales Representative"
$ws.Range("D4").Value = "Kartepe, Kocaeli, Türkiye"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "5 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "Rhino Tank"

# Row 5
$ws.Range("A5").Value = "1331094"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1331094"
$ws.Range("C5").Value = "[GBS] Transportation Sourcing Assistant"
$ws.Range("D5").Value = "Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany"
$ws.Range("E5").Value = "Yes"
$ws.Range("F5").Value = "124 applicants"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "DHL Group"

# Row 6
$ws.Range("A6").Value = "1330306"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1330306"
$ws.Range("C6").Value = "Neuro-Marketing & Communications Intern"
$ws.Range("D6").Value = "Amman, Jordan"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "56 applicants"
$ws.Range("G6").Value = "9 - 12 Weeks"
$ws.Range("H6").Value = "Amoux Group"

# Row 7 (new)
$ws.Range("A7").Value = "1327099"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1327099"
$ws.Range("C7").Value = "ACE Program | French Accounts Receivable Specialist"
$ws.Range("D7").Value = "Mumbai, Maharashtra, India"
$ws.Range("E7").Value = "Yes"
$ws.Range("F7").Value = "1 applicant"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "Tata Consultancy Services Ltd."

# Row 8 (new)
$ws.Range("A8").Value = "1311536"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1311536"
$ws.Range("C8").Value = "Accelerate Romania | Managing Co-founder"
$ws.Range("D8").Value = "Bucharest, Romania"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "179 applicants"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "Skulptor"

# Restore plain "Normal" style on the ID column so no numFmt leaks onto the
# cells themselves (values stay text via the model's stored type).
$ws.Range("A2:A8").Style = "Normal"

# Re-apply the yellow fill to the "Yes" cells since resetting A2:A8's style
# only touched column A, but make sure E-column formatting still holds.
$ws.Range("E2").Interior.Color = 65535
$ws.Range("E3").Interior.Color = 65535
$ws.Range("E5").Interior.Color = 65535
$ws.Range("E7").Interior.Color = 65535
